$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1304
$ws.Range("J17").Value = 1304
$ws.Range("L17").Value = 3912
$ws.Range("N17").Value = -4248

$ws.Range("H43").Value = 2088.3333
$ws.Range("J43").Value = 1179.8334
$ws.Range("L43").Value = 1179.8334
$ws.Range("N43").Value = -1317.8334

$ws.Range("H112").Value = 1034
$ws.Range("J112").Value = 1058.7273
$ws.Range("L112").Value = 3176.1819
$ws.Range("N112").Value = -5392.1819

$ws.Range("H129").Value = 2477.6611
$ws.Range("I129").Value = 5762.579
$ws.Range("J129").Value = 917.325
$ws.Range("K129").Value = 17287.737
$ws.Range("L129").Value = 2751.975
$ws.Range("M129").Value = -12287.737
$ws.Range("N129").Value = -12751.975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9509.666999999999
$ws.Range("I32").Value = 8450.290999999999
$ws.Range("J32").Value = 21162.8
$ws.Range("K32").Value = 8450.290999999999
$ws.Range("L32").Value = 21162.8
$ws.Range("M32").Value = -8163.290999999999
$ws.Range("N32").Value = -21736.8

$ws.Range("H39").Value = 6000
$ws.Range("I39").Value = 6000
$ws.Range("K39").Value = 6000
$ws.Range("M39").Value = -5480

$ws.Range("H40").Value = 32031
$ws.Range("J40").Value = 32031
$ws.Range("L40").Value = 32031
$ws.Range("N40").Value = -32383

$ws.Range("H46").Value = 3582.3333
$ws.Range("J46").Value = 2873.5
$ws.Range("L46").Value = 2873.5
$ws.Range("N46").Value = -3511.5

$ws.Range("H74").Value = 1457.8636
$ws.Range("I74").Value = 1670.2142
$ws.Range("J74").Value = 1086.25
$ws.Range("K74").Value = 1670.2142
$ws.Range("L74").Value = 1086.25
$ws.Range("M74").Value = -796.2141999999999
$ws.Range("N74").Value = -2834.25

$ws.Range("H77").Value = 1457.8636
$ws.Range("I77").Value = 1670.2142
$ws.Range("J77").Value = 1086.25
$ws.Range("K77").Value = 8351.071
$ws.Range("L77").Value = 5431.25
$ws.Range("M77").Value = -3983.071
$ws.Range("N77").Value = -14167.25

$ws.Range("H122").Value = 2137.9
$ws.Range("I122").Value = 2495
$ws.Range("J122").Value = 1899.8334
$ws.Range("K122").Value = 7485
$ws.Range("L122").Value = 5699.5002
$ws.Range("M122").Value = -5035
$ws.Range("N122").Value = -10599.5002

$ws.Range("H123").Value = 45001
$ws.Range("J123").Value = 45001
$ws.Range("L123").Value = 45001
$ws.Range("N123").Value = -54801

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 3914.907
$ws.Range("I132").Value = 3681.5
$ws.Range("J132").Value = 5115.2856
$ws.Range("K132").Value = 11044.5
$ws.Range("L132").Value = 15345.8568
$ws.Range("M132").Value = -8514.5
$ws.Range("N132").Value = -20405.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H134").Value = 2457.4243
$ws.Range("I134").Value = 2175.037
$ws.Range("J134").Value = 3728.1667
$ws.Range("K134").Value = 6525.110999999999
$ws.Range("L134").Value = 11184.5001
$ws.Range("M134").Value = -3990.110999999999
$ws.Range("N134").Value = -16254.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 22505.5
$ws.Range("I32").Value = 20000
$ws.Range("J32").Value = 25011
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 25011
$ws.Range("M32").Value = -19684
$ws.Range("N32").Value = -25643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 72.3
$ws.Range("I8").Value = 72.3
$ws.Range("K8").Value = 216.9
$ws.Range("M8").Value = -77.89999999999998

$ws.Range("H34").Value = 1019.38464
$ws.Range("I34").Value = 530.4
$ws.Range("J34").Value = 1325
$ws.Range("K34").Value = 1591.2
$ws.Range("L34").Value = 3975
$ws.Range("M34").Value = -1507.2
$ws.Range("N34").Value = -4143

$ws.Range("H56").Value = 3984.111
$ws.Range("I56").Value = 3984.111
$ws.Range("K56").Value = 3984.111
$ws.Range("M56").Value = -3454.111

$ws.Range("H113").Value = 777.02563
$ws.Range("J113").Value = 540.2174
$ws.Range("L113").Value = 1620.6522
$ws.Range("N113").Value = -5960.6522

$ws.Range("H121").Value = 7608.524
$ws.Range("I121").Value = 9403.799999999999
$ws.Range("J121").Value = 7047.5
$ws.Range("K121").Value = 28211.4
$ws.Range("L121").Value = 21142.5
$ws.Range("M121").Value = -26901.4
$ws.Range("N121").Value = -23762.5

$ws.Range("H122").Value = 446.9355
$ws.Range("I122").Value = 360.125
$ws.Range("J122").Value = 477.13043
$ws.Range("K122").Value = 3241.125
$ws.Range("L122").Value = 4294.17387
$ws.Range("M122").Value = -791.125
$ws.Range("N122").Value = -9194.173869999999

$ws.Range("H131").Value = 835.45
$ws.Range("I131").Value = 533
$ws.Range("J131").Value = 844.80414
$ws.Range("K131").Value = 1599
$ws.Range("L131").Value = 2534.41242
$ws.Range("M131").Value = 3441
$ws.Range("N131").Value = -12614.41242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4513.8
$ws.Range("I43").Value = 1775
$ws.Range("K43").Value = 1775
$ws.Range("M43").Value = -1624

$ws.Range("H46").Value = 12462.25
$ws.Range("J46").Value = 12462.25
$ws.Range("L46").Value = 12462.25
$ws.Range("N46").Value = -12774.25

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 100104750
$ws.Range("I80").Value = 250258750
$ws.Range("J80").Value = 2083.3333
$ws.Range("K80").Value = 250258750
$ws.Range("L80").Value = 2083.3333
$ws.Range("M80").Value = -250257752
$ws.Range("N80").Value = -4079.3333

$ws.Range("H83").Value = 100104750
$ws.Range("I83").Value = 250258750
$ws.Range("J83").Value = 2083.3333
$ws.Range("K83").Value = 1251293750
$ws.Range("L83").Value = 10416.6665
$ws.Range("M83").Value = -1251288758
$ws.Range("N83").Value = -20400.6665

$ws.Range("H102").Value = 2467.2273
$ws.Range("I102").Value = 2308.6
$ws.Range("K102").Value = 2308.6
$ws.Range("M102").Value = -686.5999999999999

$ws.Range("H139").Value = 49884
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49884
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49884
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -60164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1276.7894
$ws.Range("I22").Value = 2664.75
$ws.Range("J22").Value = 906.6667
$ws.Range("K22").Value = 2664.75
$ws.Range("L22").Value = 906.6667
$ws.Range("M22").Value = -2369.75
$ws.Range("N22").Value = -1496.6667

$ws.Range("H27").Value = 1276.7894
$ws.Range("I27").Value = 2664.75
$ws.Range("J27").Value = 906.6667
$ws.Range("K27").Value = 2664.75
$ws.Range("L27").Value = 906.6667
$ws.Range("M27").Value = -2557.75
$ws.Range("N27").Value = -1120.6667

$ws.Range("H40").Value = 61097.766
$ws.Range("I40").Value = 168333.67
$ws.Range("J40").Value = 2605.4546
$ws.Range("K40").Value = 168333.67
$ws.Range("L40").Value = 2605.4546
$ws.Range("M40").Value = -168197.67
$ws.Range("N40").Value = -2877.4546

$ws.Range("H68").Value = 4839
$ws.Range("I68").Value = 2433.3333
$ws.Range("J68").Value = 5640.8887
$ws.Range("K68").Value = 2433.3333
$ws.Range("L68").Value = 5640.8887
$ws.Range("M68").Value = -1684.3333
$ws.Range("N68").Value = -7138.8887

$ws.Range("H71").Value = 4839
$ws.Range("I71").Value = 2433.3333
$ws.Range("J71").Value = 5640.8887
$ws.Range("K71").Value = 12166.6665
$ws.Range("L71").Value = 28204.4435
$ws.Range("M71").Value = -8422.666499999999
$ws.Range("N71").Value = -35692.4435

$ws.Range("H122").Value = 3122.6667
$ws.Range("I122").Value = 3122.6667
$ws.Range("K122").Value = 9368.000100000001
$ws.Range("M122").Value = -6918.000100000001

$ws.Range("H136").Value = 1712.5927
$ws.Range("I136").Value = 1484.8096
$ws.Range("J136").Value = 2509.8333
$ws.Range("K136").Value = 4454.4288
$ws.Range("L136").Value = 7529.499899999999
$ws.Range("M136").Value = -1904.4288
$ws.Range("N136").Value = -12629.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10933.333
$ws.Range("J45").Value = 10933.333
$ws.Range("L45").Value = 10933.333
$ws.Range("N45").Value = -11915.333

$ws.Range("H136").Value = 1677.3334
$ws.Range("I136").Value = 695.1923
$ws.Range("J136").Value = 2698.76
$ws.Range("K136").Value = 2085.5769
$ws.Range("L136").Value = 8096.280000000001
$ws.Range("M136").Value = 464.4231
$ws.Range("N136").Value = -13196.28
